$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting -----------------------------------------------------------
# The header row (A1:B1) picks up the (invisible/no-op) border + explicit
# font that the data rows already carried, by copying the format down from
# row 2 onto row 1.
$ws.Range("A2:B2").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The data rows (and the still-blank row 17, which shares the same cell
# style) switch from the theme-based font color to an explicit black.
$ws.Range("A2:B17").Font.Color = 0

# Data rows grow a bit taller.
$ws.Range("2:16").RowHeight = 23.25

# --- Values -----------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 0

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 0

$ws.Range("A6").Value = 0
$ws.Range("B6").Value = 0

$ws.Range("A7").Value = 0
$ws.Range("B7").Value = 0

$ws.Range("A8").Value = 0
$ws.Range("B8").Value = 0

$ws.Range("A9").Value = 0
$ws.Range("B9").Value = 0

$ws.Range("A10").Value = 0
$ws.Range("B10").Value = 0

$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 0

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = 0

$ws.Range("A13").Value = 1
$ws.Range("B13").Value = 0

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 0
